$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.120168333333333
$ws.Range("H2").Value = 3.360505
$ws.Range("I2").Value = 0.001768092629909379
$ws.Range("J2").Value = 0.001768092629909379
$ws.Range("M2").Value = 8.142376000000001
$ws.Range("N2").Value = 24.427128
$ws.Range("O2").Value = 0.1741313933276368
$ws.Range("P2").Value = 0.1741313933276368
$ws.Range("Q2").Value = 9.120831753293333
$ws.Range("R2").Value = 82.08748577964001
$ws.Range("S2").Value = 0.0003078804331784457
$ws.Range("T2").Value = 0.0003078804331784457
$ws.Range("G3").Value = 1.120168333333333
$ws.Range("H3").Value = 3.360505
$ws.Range("I3").Value = 0.001768092629909379
$ws.Range("J3").Value = 0.001768092629909379
$ws.Range("O3").Value = 0.5205382400466131
$ws.Range("P3").Value = 0.5205382400466131
$ws.Range("Q3").Value = 27.26528294462944
$ws.Range("R3").Value = 245.387546501665
$ws.Range("S3").Value = 0.0009203598258124156
$ws.Range("T3").Value = 0.0009203598258124156
$ws.Range("G4").Value = 1.120168333333333
$ws.Range("H4").Value = 3.360505
$ws.Range("I4").Value = 0.001768092629909379
$ws.Range("J4").Value = 0.001768092629909379
$ws.Range("O4").Value = 0.3053303666257501
$ws.Range("P4").Value = 0.3053303666257501
$ws.Range("Q4").Value = 15.99290541439
$ws.Range("R4").Value = 143.93614872951
$ws.Range("S4").Value = 0.0005398523709185173
$ws.Range("T4").Value = 0.0005398523709185173
$ws.Range("I5").Value = 0.9534130698726969
$ws.Range("J5").Value = 0.9534130698726969
$ws.Range("M5").Value = 8.142376000000001
$ws.Range("N5").Value = 24.427128
$ws.Range("O5").Value = 0.1741313933276368
$ws.Range("P5").Value = 0.1741313933276368
$ws.Range("Q5").Value = 4918.249222126709
$ws.Range("R5").Value = 44264.24299914039
$ws.Range("S5").Value = 0.1660191462737122
$ws.Range("T5").Value = 0.1660191462737122
$ws.Range("I6").Value = 0.9534130698726969
$ws.Range("J6").Value = 0.9534130698726969
$ws.Range("O6").Value = 0.5205382400466131
$ws.Range("P6").Value = 0.5205382400466131
$ws.Range("S6").Value = 0.4962879614289722
$ws.Range("T6").Value = 0.4962879614289722
$ws.Range("I7").Value = 0.9534130698726969
$ws.Range("J7").Value = 0.9534130698726969
$ws.Range("O7").Value = 0.3053303666257501
$ws.Range("P7").Value = 0.3053303666257501
$ws.Range("S7").Value = 0.2911059621700125
$ws.Range("T7").Value = 0.2911059621700125
$ws.Range("H8").Value = 85.18441
$ws.Range("I8").Value = 0.04481883749739363
$ws.Range("J8").Value = 0.04481883749739363
$ws.Range("M8").Value = 8.142376000000001
$ws.Range("N8").Value = 24.427128
$ws.Range("O8").Value = 0.1741313933276368
$ws.Range("P8").Value = 0.1741313933276368
$ws.Range("Q8").Value = 231.2011651860533
$ws.Range("R8").Value = 2080.81048667448
$ws.Range("S8").Value = 0.007804366620746086
$ws.Range("T8").Value = 0.007804366620746086
$ws.Range("H9").Value = 85.18441
$ws.Range("I9").Value = 0.04481883749739363
$ws.Range("J9").Value = 0.04481883749739363
$ws.Range("O9").Value = 0.5205382400466131
$ws.Range("P9").Value = 0.5205382400466131
$ws.Range("Q9").Value = 691.1392904106144
$ws.Range("S9").Value = 0.02332991879182843
$ws.Range("T9").Value = 0.02332991879182843
$ws.Range("H10").Value = 85.18441
$ws.Range("I10").Value = 0.04481883749739363
$ws.Range("J10").Value = 0.04481883749739363
$ws.Range("O10").Value = 0.3053303666257501
$ws.Range("P10").Value = 0.3053303666257501
$ws.Range("Q10").Value = 405.39925157398
$ws.Range("R10").Value = 3648.59326416582
$ws.Range("S10").Value = 0.01368455208481911
$ws.Range("T10").Value = 0.01368455208481911
